# Reverses the order of the "Periodo Mora" rows (2203..2209 -> 2209..2203)
# in the account-statement table, keeping the "Valor Mora" amount attached
# to its original period (2209 keeps 72800, the rest keep 84000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-22, column E holds the "Periodo Mora" text; column F holds the
# "Valor Mora" amount for that same period. Re-write both columns in
# reverse row order so each period keeps its matching amount.
$periods = @("2203", "2204", "2205", "2206", "2207", "2208", "2209")
$values  = @(84000, 84000, 84000, 84000, 84000, 84000, 72800)

$startRow = 16
$rowCount = $periods.Length

for ($i = 0; $i -lt $rowCount; $i++) {
    $destRow = $startRow + ($rowCount - 1 - $i)

    $ws.Range("E$destRow").Value = $periods[$i]
    $ws.Range("F$destRow").Value = $values[$i]
}
